# Fruta / hortaliza, semanal
# Insert a new weekly data point as row 64, shifting the existing rows
# 64-74 down to 65-75 (matching the canonical OOXML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 64; all rows below (including their
# formatting) shift down by one.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new observation.
$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value = "Los Lagos"
$ws.Cells.Item(64, 4).Value = 44504
$ws.Cells.Item(64, 5).Value = 10
$ws.Cells.Item(64, 6).Value = 100112052
$ws.Cells.Item(64, 7).Value = "Albahaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 60
$ws.Cells.Item(64, 11).Value = 7000
$ws.Cells.Item(64, 12).Value = 7000
$ws.Cells.Item(64, 13).Value = 7000
$ws.Cells.Item(64, 14).Value = "`$/paquete"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 7000
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
